# feat(phone): add stylesheet 720x1280
# Update the "O" column (720x1280 stylesheet) values and move the
# worksheet selection to K7.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update input values driving the dependent formulas in column O.
$ws.Range("O13").Value = 80
$ws.Range("O14").Value = 28

# Recalculate so dependent formula cells (O15, O16, O22, ...) refresh.
$excel.Calculate()

# Move the active selection to K7, matching the saved view state.
$ws.Range("K7").Select()
